$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# Fix a typo in a stock name: "中鋼 •" -> "中鋼•" (drop the stray space before the bullet)
$ws.Range("B9").Value = "中鋼•"

# Correct the total value for the last stock row (聯格科技): "1,000,000" -> "1000000"
# (force text so it keeps matching the string-typed cell, not a number)
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "1000000"

# Insert a new column H ("property_category") before the existing date/legislator columns,
# shifting date -> I, legislator_name -> J, legislator_id -> K
$ws.Columns.Item(8).Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"
$ws.Range("H5").Value = "stock"
$ws.Range("H6").Value = "stock"
$ws.Range("H7").Value = "stock"
$ws.Range("H8").Value = "stock"
$ws.Range("H9").Value = "stock"
$ws.Range("H10").Value = "stock"
$ws.Range("H11").Value = "stock"
$ws.Range("H12").Value = "stock"
$ws.Range("H13").Value = "stock"
$ws.Range("H14").Value = "stock"
$ws.Range("H15").Value = "stock"
